$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Insert a new column before column D (so existing D..F shift to E..G)
$ws.Columns("D").Insert()

# Match the width of the neighboring "execute priority" column (C)
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Set the header for the new "version" column
$ws.Range("D1").Value = "version"

# Fill the new column with "N/A" for all data rows
$ws.Range("D2:D7").Value = "N/A"

# Update the selection to match the target state
$ws.Range("D2:D7").Select()
